$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N2").ClearContents()
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -162.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 275.5
$ws.Range("H2").Value = 275.5
$ws.Range("I2").Value = 275.5
$ws.Range("N9").Value = -1138
$ws.Range("L9").Value = 800
$ws.Range("K9").Value = 1181.375
$ws.Range("I9").Value = 1181.375
$ws.Range("J9").Value = 800
$ws.Range("H9").Value = 1139
$ws.Range("M9").Value = -1012.375
$ws.Range("N54").ClearContents()
$ws.Range("H54").Value = 31000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("K100").Value = 1871.2
$ws.Range("M100").Value = -1330.2
$ws.Range("I100").Value = 1871.2
$ws.Range("H100").Value = 5246.273
$ws.Range("M132").Value = -1095.3125
$ws.Range("L132").Value = 16900.0005
$ws.Range("K132").Value = 3625.3125
$ws.Range("I132").Value = 1208.4375
$ws.Range("H132").Value = 1907.1052
$ws.Range("N132").Value = -21960.0005
$ws.Range("J132").Value = 5633.3335
$ws.Range("J133").Value = 59616.848
$ws.Range("H133").Value = 59616.848
$ws.Range("N133").Value = -69736.848
$ws.Range("L133").Value = 59616.848
$ws.Range("K137").Value = 9351.1428
$ws.Range("I137").Value = 3117.0476
$ws.Range("H137").Value = 4681.0225
$ws.Range("M137").Value = -6801.1428

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 9514.833000000001
$ws.Range("L32").Value = 9514.833000000001
$ws.Range("H32").Value = 4840.186
$ws.Range("N32").Value = -10088.833
$ws.Range("M45").Value = -4843.6
$ws.Range("H45").Value = 9091.75
$ws.Range("K45").Value = 5220.6
$ws.Range("I45").Value = 5220.6
$ws.Range("H63").Value = 2773.7144
$ws.Range("I63").Value = 2773.7144
$ws.Range("M63").Value = -2087.7144
$ws.Range("K63").Value = 2773.7144
$ws.Range("I66").Value = 2773.7144
$ws.Range("H66").Value = 2773.7144
$ws.Range("M66").Value = -10436.572
$ws.Range("K66").Value = 13868.572
$ws.Range("H110").Value = 209346.3
$ws.Range("K110").Value = 264177.38
$ws.Range("I110").Value = 264177.38
$ws.Range("M110").Value = -262132.38
$ws.Range("M132").Value = -2100.1112
$ws.Range("L132").Value = 27637.089
$ws.Range("K132").Value = 4630.1112
$ws.Range("I132").Value = 1543.3704
$ws.Range("H132").Value = 4986.592
$ws.Range("N132").Value = -32697.089
$ws.Range("J132").Value = 9212.362999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N94").ClearContents()
$ws.Range("H94").Value = 1518.5
$ws.Range("J94").Value = 0
$ws.Range("I94").Value = 1518.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1067.5
$ws.Range("K94").Value = 1518.5
$ws.Range("L109").Value = 60000
$ws.Range("J109").Value = 60000
$ws.Range("N109").Value = -62774
$ws.Range("H109").Value = 60000
$ws.Range("H134").Value = 4122.9487
$ws.Range("I134").Value = 2744.5
$ws.Range("K134").Value = 8233.5
$ws.Range("M134").Value = -5698.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J58").Value = 5662.85
$ws.Range("N58").Value = -6068.85
$ws.Range("M58").Value = -556676.25
$ws.Range("I58").Value = 556879.25
$ws.Range("L58").Value = 5662.85
$ws.Range("K58").Value = 556879.25
$ws.Range("H58").Value = 266765.34
$ws.Range("M132").Value = -5189.500100000001
$ws.Range("L132").Value = 17358.2139
$ws.Range("K132").Value = 7719.500100000001
$ws.Range("I132").Value = 2573.1667
$ws.Range("H132").Value = 3978.8125
$ws.Range("N132").Value = -22418.2139
$ws.Range("J132").Value = 5786.0713
$ws.Range("H134").Value = 5573.409
$ws.Range("I134").Value = 4741.3335
$ws.Range("K134").Value = 14224.0005
$ws.Range("M134").Value = -11689.0005
$ws.Range("M136").Value = -1668087.75
$ws.Range("N136").Value = -22088.55
$ws.Range("J136").Value = 5662.85
$ws.Range("I136").Value = 556879.25
$ws.Range("L136").Value = 16988.55
$ws.Range("H136").Value = 266765.34
$ws.Range("K136").Value = 1670637.75
$ws.Range("L138").Value = 49990
$ws.Range("H138").Value = 49990
$ws.Range("N138").Value = -60270
$ws.Range("J138").Value = 49990

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L10").Value = 30
$ws.Range("N10").Value = -308
$ws.Range("H10").Value = 53.4
$ws.Range("J10").Value = 10
$ws.Range("N11").Value = -1165
$ws.Range("M11").Value = -6667459.600000001
$ws.Range("K11").Value = 6667599.600000001
$ws.Range("I11").Value = 2222533.2
$ws.Range("L11").Value = 885
$ws.Range("J11").Value = 295
$ws.Range("H11").Value = 2000309.4
$ws.Range("H12").Value = 316.08334
$ws.Range("K12").Value = 148.5
$ws.Range("M12").Value = 24.5
$ws.Range("I12").Value = 49.5
$ws.Range("N101").Value = -44502.8
$ws.Range("J101").Value = 13211.6
$ws.Range("H101").Value = 13211.6
$ws.Range("L101").Value = 39634.8
$ws.Range("I121").Value = 1280.5714
$ws.Range("H121").Value = 557206.75
$ws.Range("M121").Value = -2531.7142
$ws.Range("K121").Value = 3841.7142
$ws.Range("J129").Value = 202973
$ws.Range("H129").Value = 63809.25
$ws.Range("L129").Value = 608919
$ws.Range("N129").Value = -618919

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M2").Value = -99.19999999999999
$ws.Range("K2").Value = 212.2
$ws.Range("H2").Value = 165.53847
$ws.Range("I2").Value = 212.2
$ws.Range("H33").Value = 14900
$ws.Range("J33").Value = 14900
$ws.Range("N33").Value = -15404
$ws.Range("L33").Value = 14900
$ws.Range("L34").Value = 47025.5
$ws.Range("H34").Value = 48016
$ws.Range("J34").Value = 47025.5
$ws.Range("N34").Value = -47561.5
$ws.Range("N40").Value = -5302
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N74").Value = -31871
$ws.Range("H74").Value = 29999
$ws.Range("L74").Value = 29999
$ws.Range("J74").Value = 29999
$ws.Range("H76").Value = 48016
$ws.Range("L76").Value = 47025.5
$ws.Range("N76").Value = -47655.5
$ws.Range("J76").Value = 47025.5
$ws.Range("N77").Value = -99357
$ws.Range("H77").Value = 29999
$ws.Range("L77").Value = 89997
$ws.Range("J77").Value = 29999
$ws.Range("J79").Value = 47025.5
$ws.Range("L79").Value = 47025.5
$ws.Range("H79").Value = 48016
$ws.Range("N79").Value = -49209.5
$ws.Range("M132").Value = -3344160.8
$ws.Range("L132").Value = 20248.875
$ws.Range("K132").Value = 3346690.8
$ws.Range("I132").Value = 1115563.6
$ws.Range("H132").Value = 593768.8
$ws.Range("N132").Value = -25308.875
$ws.Range("J132").Value = 6749.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N16").ClearContents()
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("J126").Value = 4566
$ws.Range("N126").Value = -18638
$ws.Range("H126").Value = 3182
$ws.Range("L126").Value = 13698
$ws.Range("M136").Value = -6165.1875
$ws.Range("I136").Value = 2905.0625
$ws.Range("H136").Value = 2850.7058
$ws.Range("K136").Value = 2905.0625
